$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Variables")
$ws.Rows.Item(5).Insert()
$ws.Range("A5").Value = "age_trimester"
$ws.Range("B5").Value = "integer"
$ws.Range("C5").Value = "numeric"
$ws.Range("D5").Value = "Age of the child in trimesters"
